$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the last week sheet ("Nädal 9") to create the new week
#    ("Nädal 10"), placed right after it.
# ---------------------------------------------------------------------
$week9 = $wb.Worksheets.Item("Nädal 9")
$week9.Copy([System.Reflection.Missing]::Value, $week9)
$week10 = $wb.Worksheets.Item($week9.Index + 1)
$week10.Name = "Nädal 10"

# ---------------------------------------------------------------------
# 2. Finish filling in week 9's log (last row: Friday 03.04.2020,
#    kontrolltöö prep video).
# ---------------------------------------------------------------------
$week9.Range("B13").Value2 = 43924
$week9.Range("C13").Value2 = 0.90277777777777779
$week9.Range("D13").Value2 = 0
$week9.Range("F13").Value2 = 140
$week9.Range("G13").Value2 = "video"
$week9.Range("H13").Value2 = "V48"

# Week 9 is now a completed week: tab turns pink/red and it is no longer
# the highlighted/active sheet.
$week9.Tab.Color = 8420607
$week9.Range("A1:J17").Select()

# ---------------------------------------------------------------------
# 3. Reset week 10's data area - it was copied verbatim from week 9, so
#    clear out all the old entries before putting in the new week's data.
# ---------------------------------------------------------------------
$week10.Range("B4:J12").ClearContents()

# Header date range for the new week
$week10.Range("G1").Value2 = "03.04.2020 - 09.04.2020"

# Row 1 of the log: Saturday 04.04.2020, the exam (kontrolltöö)
$week10.Range("B4").Value2 = 43925
$week10.Range("C4").Value2 = 0.50694444444444442
$week10.Range("D4").Value2 = 0.54861111111111105
$week10.Range("F4").Value2 = 60
$week10.Range("G4").Value2 = "kontrolltöö"
$week10.Range("H4").Value2 = "KT, tulemus 9,5/10"
$week10.Range("I4").Value2 = "x"

# Row 2 of the log: still same day, watching video about fasaadi tests
$week10.Range("C5").Value2 = 0.90625
$week10.Range("G5").Value2 = "video"
$week10.Range("H5").Value2 = "V48, 49"

# Week 10 becomes the active/highlighted (green) tab with the cursor on G6
$week10.Range("G6").Select()

Write-Host "Edit complete"
